$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated election result figures for row 2 (SETÚBAL / ALCÁCER DO SAL)
$ws.Range("H2").Value  = 91
$ws.Range("I2").Value  = 282
$ws.Range("J2").Value  = 1219
$ws.Range("K2").Value  = 4
$ws.Range("L2").Value  = 330
$ws.Range("M2").Value  = 25
$ws.Range("N2").Value  = 207
$ws.Range("O2").Value  = 1
$ws.Range("P2").Value  = 8
$ws.Range("Q2").Value  = 1
$ws.Range("R2").Value  = 12
$ws.Range("S2").Value  = 133
$ws.Range("T2").Value  = 198
$ws.Range("U2").Value  = 17
$ws.Range("V2").Value  = 1907
$ws.Range("W2").Value  = 1
$ws.Range("X2").Value  = 1939
$ws.Range("Y2").Value  = 2
$ws.Range("Z2").Value  = 33
$ws.Range("AA2").Value = 13

$wb.Save()
